# Update the "USD Amount" figure for the Deposit/Crypto/Roobic row (row 2)
# on SheetName1 from 28260 to 29439.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")
$ws.Range("T2").Value = 29439
